$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 226, shifting existing rows 226:323 down to 227:324
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new weekly data point
$ws.Cells.Item(226, 1).Value = 10
$ws.Cells.Item(226, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(226, 3).Value = "La Araucanía"
$ws.Cells.Item(226, 4).Value = 45141
$ws.Cells.Item(226, 5).Value = 9
$ws.Cells.Item(226, 6).Value = 100112005
$ws.Cells.Item(226, 7).Value = "Puerro"
$ws.Cells.Item(226, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 120
$ws.Cells.Item(226, 11).Value = 8000
$ws.Cells.Item(226, 12).Value = 8000
$ws.Cells.Item(226, 13).Value = 8000
$ws.Cells.Item(226, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(226, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(226, 16).Value = 667
$ws.Cells.Item(226, 17).Value = 12
$ws.Cells.Item(226, 18).Value = "Hortaliza"
